$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.603458523750305
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 3.238240242004395
$ws.Range("D1").Value = 1.260423898696899
$ws.Range("E1").Value = 0.8281568884849548
